$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp text in F1
$ws.Range("F1").Value = "Last status check on: 25.01.2022 06:15"

# D4: change from text "+0.4" to numeric 0.4
$ws.Range("D4").Value = 0.4

# E4: change from text "2022-01-25 06:00:12" to a real date/time value,
# matching the numeric format used by the other Old Datum cells (E2, E3, etc.)
$ws.Range("E4").Value = 44586.25013888889
$ws.Range("E4").NumberFormat = $ws.Range("E3").NumberFormat
